$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 1) Back up the exact cell formatting that would otherwise be
#    clobbered by re-merging (Excel copies the anchor cell's format
#    onto every cell of a range when it is (re)merged). We stash a
#    copy of each relevant format in unused scratch cells far away
#    from the sheet's real data (row 100) so we can restore it again
#    once the merges have been recreated in the desired order.
# ---------------------------------------------------------------------
$ws.Range("Z100").Copy()                      # pristine / default format
$ws.Range("A100").PasteSpecial($fmt)

$ws.Range("B9").Copy()                        # row 9 interior (B9:D9)
$ws.Range("B100").PasteSpecial($fmt)

$ws.Range("B10").Copy()                       # row 10 interior (B10:D10)
$ws.Range("C100").PasteSpecial($fmt)

$ws.Range("B11").Copy()                       # row 11 interior (B11:D11)
$ws.Range("D100").PasteSpecial($fmt)

$ws.Range("C13").Copy()                       # C13
$ws.Range("E100").PasteSpecial($fmt)

$ws.Range("C14").Copy()                       # C14
$ws.Range("F100").PasteSpecial($fmt)

$ws.Range("C15").Copy()                       # C15
$ws.Range("G100").PasteSpecial($fmt)

$ws.Range("C16").Copy()                       # C16
$ws.Range("H100").PasteSpecial($fmt)

$ws.Range("B31").Copy()                       # row 31 interior (B31:C31)
$ws.Range("I100").PasteSpecial($fmt)

$ws.Range("D31").Copy()                       # D31
$ws.Range("J100").PasteSpecial($fmt)

# ---------------------------------------------------------------------
# 2) Recreate every merged range so that, once finished, the
#    <mergeCells> list is emitted in the desired final order (new
#    merges are appended at save time in the order Merge() is called).
# ---------------------------------------------------------------------
$ws.Range("A10:D10").UnMerge()
$ws.Range("B15:C15").UnMerge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("B13:C13").UnMerge()
$ws.Range("A1:D1").UnMerge()
$ws.Range("B16:C16").UnMerge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("A31:D31").UnMerge()
$ws.Range("B14:C14").UnMerge()

$ws.Range("B13:C13").Merge()
$ws.Range("A1:D1").Merge()
$ws.Range("A9:D9").Merge()
$ws.Range("B16:C16").Merge()
$ws.Range("A31:D31").Merge()
$ws.Range("B15:C15").Merge()
$ws.Range("B14:C14").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("A10:D10").Merge()

# ---------------------------------------------------------------------
# 3) Restore the original per-cell formatting that Merge() overwrote.
# ---------------------------------------------------------------------
$ws.Range("A100").Copy()
$ws.Range("B1:D1").PasteSpecial($fmt)

$ws.Range("B100").Copy()
$ws.Range("B9:D9").PasteSpecial($fmt)

$ws.Range("C100").Copy()
$ws.Range("B10:D10").PasteSpecial($fmt)

$ws.Range("D100").Copy()
$ws.Range("B11:D11").PasteSpecial($fmt)

$ws.Range("E100").Copy()
$ws.Range("C13").PasteSpecial($fmt)

$ws.Range("F100").Copy()
$ws.Range("C14").PasteSpecial($fmt)

$ws.Range("G100").Copy()
$ws.Range("C15").PasteSpecial($fmt)

$ws.Range("H100").Copy()
$ws.Range("C16").PasteSpecial($fmt)

$ws.Range("I100").Copy()
$ws.Range("B31:C31").PasteSpecial($fmt)

$ws.Range("J100").Copy()
$ws.Range("D31").PasteSpecial($fmt)

# ---------------------------------------------------------------------
# 4) Clean up the scratch cells used for the backup/restore dance.
# ---------------------------------------------------------------------
$ws.Range("A100:J100").Clear()

# ---------------------------------------------------------------------
# 5) Apply the actual data changes described by the diff.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = 45311          # date bumped by one day

$ws.Range("D14").Value = 707.256
$ws.Range("D15").Value = 863.588
$ws.Range("D16").Value = 1014.597
